$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 45029
$ws.Range("Q2").Value = "`$/bandeja 18 kilos granel"
$ws.Range("S2").Value = 528
$ws.Range("T2").Value = 18

# Row 3
$ws.Range("D3").Value = 44776
$ws.Range("L3").Value = "Primera"
$ws.Range("N3").Value = 10000
$ws.Range("O3").Value = 10000
$ws.Range("P3").Value = 10000
$ws.Range("Q3").Value = "`$/bandeja 18 kilos granel"
$ws.Range("S3").Value = 556
$ws.Range("T3").Value = 18

# Row 4
$ws.Range("D4").Value = 44776
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 8000
$ws.Range("O4").Value = 8000
$ws.Range("P4").Value = 8000
$ws.Range("S4").Value = 444

# Row 5
$ws.Range("D5").Value = 45041
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 11000
$ws.Range("O5").Value = 12000
$ws.Range("P5").Value = 11500
$ws.Range("Q5").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 639

# Row 6
$ws.Range("D6").Value = 45034
$ws.Range("M6").Value = 220
$ws.Range("N6").Value = 8500
$ws.Range("O6").Value = 9000
$ws.Range("P6").Value = 8727
$ws.Range("Q6").Value = "`$/caja 18 kilos granel"
$ws.Range("S6").Value = 485

# Row 7
$ws.Range("D7").Value = 45027
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 9000
$ws.Range("O7").Value = 10000
$ws.Range("P7").Value = 9500
$ws.Range("S7").Value = 528

# Row 8
$ws.Range("D8").Value = 44363
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 9000
$ws.Range("O8").Value = 10000
$ws.Range("P8").Value = 9500
$ws.Range("Q8").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R8").Value = "Región de O'Higgins"
$ws.Range("S8").Value = 633
$ws.Range("T8").Value = 15

# Row 9
$ws.Range("D9").Value = 44358
$ws.Range("N9").Value = 11000
$ws.Range("O9").Value = 12000
$ws.Range("P9").Value = 11500
$ws.Range("Q9").Value = "`$/caja 18 kilos granel"
$ws.Range("S9").Value = 639

# Row 10
$ws.Range("D10").Value = 44307
$ws.Range("Q10").Value = "`$/bandeja 18 kilos granel"

# Row 11
$ws.Range("D11").Value = 44307
$ws.Range("L11").Value = "Segunda"
$ws.Range("M11").Value = 50
$ws.Range("N11").Value = 8000
$ws.Range("O11").Value = 8000
$ws.Range("P11").Value = 8000
$ws.Range("Q11").Value = "`$/bandeja 18 kilos granel"
$ws.Range("S11").Value = 444

# Row 12
$ws.Range("D12").Value = 45014

# Row 13
$ws.Range("D13").Value = 44698
$ws.Range("M13").Value = 50
$ws.Range("N13").Value = 10000
$ws.Range("P13").Value = 10000
$ws.Range("S13").Value = 556

# Row 14
$ws.Range("D14").Value = 44316
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = 9000
$ws.Range("O14").Value = 10000
$ws.Range("P14").Value = 9500
$ws.Range("S14").Value = 528

# Row 15
$ws.Range("D15").Value = 44299
$ws.Range("N15").Value = 10000
$ws.Range("O15").Value = 11000
$ws.Range("P15").Value = 10500
$ws.Range("Q15").Value = "`$/caja 18 kilos granel"
$ws.Range("R15").Value = "Región del Maule"
$ws.Range("S15").Value = 583

# Row 16
$ws.Range("L16").Value = "Segunda"
$ws.Range("M16").Value = 50
$ws.Range("N16").Value = 9000
$ws.Range("O16").Value = 9000
$ws.Range("P16").Value = 9000
$ws.Range("S16").Value = 500

# Row 17
$ws.Range("D17").Value = 45050
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 140
$ws.Range("N17").Value = 11000
$ws.Range("O17").Value = 12000
$ws.Range("P17").Value = 11429
$ws.Range("Q17").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R17").Value = "Región de O'Higgins"
$ws.Range("S17").Value = 635

# Row 18
$ws.Range("D18").Value = 45013
$ws.Range("Q18").Value = "`$/bandeja 18 kilos granel"
$ws.Range("S18").Value = 528
$ws.Range("T18").Value = 18

# Row 19
$ws.Range("D19").Value = 44999
$ws.Range("M19").Value = 100
$ws.Range("N19").Value = 12000
$ws.Range("P19").Value = 12000
$ws.Range("Q19").Value = "`$/bandeja 18 kilos granel"
$ws.Range("S19").Value = 667

# Row 20
$ws.Range("D20").Value = 44999
$ws.Range("L20").Value = "Segunda"
$ws.Range("M20").Value = 100

# Row 21
$ws.Range("D21").Value = 44272
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 100
$ws.Range("N21").Value = 9000
$ws.Range("O21").Value = 10000
$ws.Range("P21").Value = 9500
$ws.Range("Q21").Value = "`$/caja 15 kilos granel"
$ws.Range("S21").Value = 633
$ws.Range("T21").Value = 15

# Row 22
$ws.Range("D22").Value = 44272
$ws.Range("L22").Value = "Segunda"
$ws.Range("M22").Value = 50
$ws.Range("N22").Value = 8000
$ws.Range("O22").Value = 8000
$ws.Range("P22").Value = 8000
$ws.Range("Q22").Value = "`$/caja 15 kilos granel"
$ws.Range("S22").Value = 533
$ws.Range("T22").Value = 15

# Row 23
$ws.Range("D23").Value = 45076
$ws.Range("M23").Value = 150
$ws.Range("N23").Value = 10000
$ws.Range("O23").Value = 11000
$ws.Range("P23").Value = 10467
$ws.Range("Q23").Value = "`$/caja 18 kilos granel"
$ws.Range("R23").Value = "Provincia de Curicó"
$ws.Range("S23").Value = 582

# Row 24
$ws.Range("D24").Value = 44425
$ws.Range("O24").Value = 13000
$ws.Range("P24").Value = 12500
$ws.Range("S24").Value = 694

# Row 25
$ws.Range("D25").Value = 45037
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 250
$ws.Range("N25").Value = 9000
$ws.Range("O25").Value = 9500
$ws.Range("P25").Value = 9200
$ws.Range("Q25").Value = "`$/caja 18 kilos granel"
$ws.Range("R25").Value = "Provincia de Curicó"
$ws.Range("S25").Value = 511
